# "them code lay data tu plc" - update the PLC-data-type rows on Sheet1:
# several "Real" entries become "Integer", the two now-unused 0.5 / 0.6
# values on rows 2-3 collapse to "0", the workbook window position moves,
# and the active selection shifts to F6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column C ("Kieu"/Type) switches from "Real" to "Integer" for the first
# two machine blocks' header rows (2,3) and the start-of-block rows that
# follow every 20 rows through row 183.
$typeRows = @(2,3,22,23,42,43,62,63,82,83,102,103,122,123,142,143,162,163,182,183)
foreach ($r in $typeRows) {
    $ws.Cells.Item($r, 3).Value2 = "Integer"
}

# The first row of each block (2,22,42,...,182) also loses its top border
# (was style 7 with a top rule; becomes borderless like the row below it).
$firstTypeRows = @(2,22,42,62,82,102,122,142,162,182)
foreach ($r in $firstTypeRows) {
    $ws.Cells.Item($r, 3).Borders.Item(8).LineStyle = -4142
}

# Column G ("Gia tri"/Value) for rows 2 and 3 used to hold the now-removed
# shared strings "0.5" / "0.6" - they become "0" (shared with other rows).
$ws.Cells.Item(2, 7).Value2 = "0"
$ws.Cells.Item(3, 7).Value2 = "0"

# Window was scrolled/moved and the live selection moved from G3 to F6.
$excel.ActiveWindow.Left = 1860
$excel.ActiveWindow.Top = 0
$ws.Range("F6").Select() | Out-Null
